$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-28 17:48:35"
$ws.Range("K2").Value = "11.8 MJ/m2"
$ws.Range("E3").Value = "2026-02-28 17:48:37"
$ws.Range("O3").Value = "-1.0 °C"
$ws.Range("E4").Value = "2026-02-28 17:48:40"
$ws.Range("H4").Value = "81%"
$ws.Range("K4").Value = "6.1 MJ/m2"
$ws.Range("O4").Value = "11.2 °C"
$ws.Range("E5").Value = "2026-02-28 17:48:42"
$ws.Range("I5").Value = "0.8 mm"
$ws.Range("N5").Value = "-2.4 °C 17:20 TU"
$ws.Range("O5").Value = "-1.1 °C"
$ws.Range("E6").Value = "2026-02-28 17:48:45"
$ws.Range("E7").Value = "2026-02-28 17:48:47"
$ws.Range("O7").Value = "13.2 °C"
$ws.Range("E8").Value = "2026-02-28 17:48:50"
$ws.Range("O8").Value = "9.9 °C"
$ws.Range("E9").Value = "2026-02-28 17:48:52"
$ws.Range("E10").Value = "2026-02-28 17:48:55"
$ws.Range("E11").Value = "2026-02-28 17:48:57"
$ws.Range("O11").Value = "7.0 °C"
$ws.Range("E12").Value = "2026-02-28 17:49:00"
$ws.Range("H12").Value = "84%"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-28 17:49:02"
$ws.Range("J13").Value = "1024.0 hPa"
$ws.Range("K13").Value = "12.0 MJ/m2"
$ws.Range("O13").Value = "6.4 °C"
$ws.Range("E14").Value = "2026-02-28 17:49:04"
$ws.Range("E15").Value = "2026-02-28 17:49:07"
$ws.Range("O15").Value = "11.3 °C"
$ws.Range("E16").Value = "2026-02-28 17:49:09"
$ws.Range("E17").Value = "2026-02-28 17:49:12"
$ws.Range("E18").Value = "2026-02-28 17:49:14"
$ws.Range("H18").Value = "82%"
$ws.Range("O18").Value = "11.9 °C"
$ws.Range("E19").Value = "2026-02-28 17:49:17"
$ws.Range("H19").Value = "76%"
$ws.Range("K19").Value = "7.8 MJ/m2"
$ws.Range("E20").Value = "2026-02-28 17:49:20"
$ws.Range("H20").Value = "57%"
$ws.Range("K20").Value = "14.0 MJ/m2"
$ws.Range("E21").Value = "2026-02-28 17:49:22"
$ws.Range("K21").Value = "11.0 MJ/m2"
$ws.Range("O21").Value = "7.5 °C"
$ws.Range("E22").Value = "2026-02-28 17:49:25"
$ws.Range("H22").Value = "64%"
$ws.Range("K22").Value = "11.5 MJ/m2"
$ws.Range("E23").Value = "2026-02-28 17:49:27"
$ws.Range("H23").Value = "67%"
$ws.Range("N23").Value = "-1.3 °C 17:29 TU"
$ws.Range("E24").Value = "2026-02-28 17:49:30"
$ws.Range("O24").Value = "8.3 °C"
$ws.Range("E25").Value = "2026-02-28 17:49:32"
$ws.Range("H25").Value = "57%"
$ws.Range("E26").Value = "2026-02-28 17:49:35"
$ws.Range("H26").Value = "78%"
$ws.Range("J26").Value = "1024.2 hPa"
$ws.Range("O26").Value = "5.0 °C"
$ws.Range("E27").Value = "2026-02-28 17:49:38"
$ws.Range("H27").Value = "49%"
$ws.Range("K27").Value = "13.9 MJ/m2"
$ws.Range("O27").Value = "2.2 °C"
$ws.Range("E28").Value = "2026-02-28 17:49:40"
$ws.Range("H28").Value = "82%"
$ws.Range("O28").Value = "9.5 °C"
$ws.Range("E29").Value = "2026-02-28 17:49:43"
$ws.Range("K29").Value = "12.4 MJ/m2"
$ws.Range("O29").Value = "11.8 °C"
$ws.Range("E30").Value = "2026-02-28 17:49:46"
$ws.Range("O30").Value = "11.0 °C"
$ws.Range("E31").Value = "2026-02-28 17:49:48"
$ws.Range("H31").Value = "83%"
$ws.Range("K31").Value = "12.0 MJ/m2"
$ws.Range("E32").Value = "2026-02-28 17:49:51"
$ws.Range("H32").Value = "97%"
$ws.Range("E33").Value = "2026-02-28 17:49:54"
$ws.Range("E34").Value = "2026-02-28 17:49:57"
$ws.Range("H34").Value = "63%"
$ws.Range("E35").Value = "2026-02-28 17:49:59"
$ws.Range("E36").Value = "2026-02-28 17:50:02"
$ws.Range("H36").Value = "82%"
$ws.Range("E37").Value = "2026-02-28 17:50:04"
$ws.Range("J37").Value = "1025.8 hPa"
$ws.Range("O37").Value = "7.0 °C"
$ws.Range("E38").Value = "2026-02-28 17:50:07"
$ws.Range("E39").Value = "2026-02-28 17:50:10"
$ws.Range("H39").Value = "58%"
$ws.Range("N39").Value = "-1.8 °C 17:20 TU"
$ws.Range("O39").Value = "-0.4 °C"
$ws.Range("E40").Value = "2026-02-28 17:50:12"
$ws.Range("H40").Value = "76%"
$ws.Range("O40").Value = "7.3 °C"
$ws.Range("E41").Value = "2026-02-28 17:50:15"
$ws.Range("E42").Value = "2026-02-28 17:50:17"
$ws.Range("E43").Value = "2026-02-28 17:50:20"
$ws.Range("O43").Value = "7.3 °C"
$ws.Range("E44").Value = "2026-02-28 17:50:22"
$ws.Range("I44").Value = "0.3 mm"
$ws.Range("O44").Value = "-0.9 °C"
$ws.Range("E45").Value = "2026-02-28 17:50:25"
$ws.Range("J45").Value = "1025.2 hPa"
$ws.Range("O45").Value = "6.5 °C"
$ws.Range("E46").Value = "2026-02-28 17:50:28"
$ws.Range("K46").Value = "4.8 MJ/m2"
$ws.Range("O46").Value = "11.7 °C"
